$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 373.8889
$ws.Range("I32").Value = 384.2
$ws.Range("J32").Value = 361
$ws.Range("K32").Value = 384.2
$ws.Range("L32").Value = 361
$ws.Range("M32").Value = -58.19999999999999
$ws.Range("N32").Value = -1013
$ws.Range("H40").Value = 1323.1538
$ws.Range("I40").Value = 1380.2
$ws.Range("J40").Value = 1287.5
$ws.Range("K40").Value = 1380.2
$ws.Range("L40").Value = 1287.5
$ws.Range("M40").Value = -1205.2
$ws.Range("N40").Value = -1637.5
$ws.Range("H55").Value = 163.25
$ws.Range("J55").Value = 214.28572
$ws.Range("L55").Value = 214.28572
$ws.Range("N55").Value = -642.28572
$ws.Range("H137").Value = 2149.5957
$ws.Range("I137").Value = 1461.3422
$ws.Range("J137").Value = 5055.5557
$ws.Range("K137").Value = 4384.0266
$ws.Range("L137").Value = 15166.6671
$ws.Range("M137").Value = -1834.0266
$ws.Range("N137").Value = -20266.6671

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 630.7179599999999
$ws.Range("I2").Value = 584.59375
$ws.Range("J2").Value = 841.5714
$ws.Range("K2").Value = 584.59375
$ws.Range("L2").Value = 841.5714
$ws.Range("M2").Value = -471.59375
$ws.Range("N2").Value = -1067.5714
$ws.Range("H61").Value = 1447.3684
$ws.Range("I61").Value = 1226.0883
$ws.Range("K61").Value = 1226.0883
$ws.Range("M61").Value = -1014.0883
$ws.Range("H74").Value = 1459.0264
$ws.Range("I74").Value = 1052.6364
$ws.Range("J74").Value = 4141.2
$ws.Range("K74").Value = 1052.6364
$ws.Range("L74").Value = 4141.2
$ws.Range("M74").Value = -178.6364000000001
$ws.Range("N74").Value = -5889.2
$ws.Range("H77").Value = 1459.0264
$ws.Range("I77").Value = 1052.6364
$ws.Range("J77").Value = 4141.2
$ws.Range("K77").Value = 5263.182000000001
$ws.Range("L77").Value = 20706
$ws.Range("M77").Value = -895.1820000000007
$ws.Range("N77").Value = -29442
$ws.Range("H116").Value = 630.7179599999999
$ws.Range("I116").Value = 584.59375
$ws.Range("J116").Value = 841.5714
$ws.Range("K116").Value = 584.59375
$ws.Range("L116").Value = 841.5714
$ws.Range("M116").Value = 1709.40625
$ws.Range("N116").Value = -5429.5714
$ws.Range("H132").Value = 3864.9722
$ws.Range("I132").Value = 3183.24
$ws.Range("K132").Value = 9549.719999999999
$ws.Range("M132").Value = -7019.719999999999
$ws.Range("H136").Value = 1447.3684
$ws.Range("I136").Value = 1226.0883
$ws.Range("K136").Value = 3678.2649
$ws.Range("M136").Value = -1128.2649

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 630.7179599999999
$ws.Range("I3").Value = 584.59375
$ws.Range("J3").Value = 841.5714
$ws.Range("K3").Value = 584.59375
$ws.Range("L3").Value = 841.5714
$ws.Range("M3").Value = -470.59375
$ws.Range("N3").Value = -1069.5714
$ws.Range("H94").Value = 1418.4615
$ws.Range("I94").Value = 1390.8334
$ws.Range("K94").Value = 1390.8334
$ws.Range("M94").Value = -939.8334
$ws.Range("H134").Value = 2608.5667
$ws.Range("I134").Value = 1575.6327
$ws.Range("K134").Value = 4726.8981
$ws.Range("M134").Value = -2191.8981

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 36000
$ws.Range("J124").Value = 36000
$ws.Range("L124").Value = 36000
$ws.Range("N124").Value = -40910

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1046.1628
$ws.Range("J5").Value = 2966.7273
$ws.Range("L5").Value = 8900.1819
$ws.Range("N5").Value = -9124.1819
$ws.Range("H17").Value = 1460
$ws.Range("I17").Value = 566.6667
$ws.Range("J17").Value = 2800
$ws.Range("K17").Value = 1700.0001
$ws.Range("L17").Value = 8400
$ws.Range("M17").Value = -1531.0001
$ws.Range("N17").Value = -8738
$ws.Range("H34").Value = 21638.6
$ws.Range("I34").Value = 36849.332
$ws.Range("J34").Value = 15119.714
$ws.Range("K34").Value = 110547.996
$ws.Range("L34").Value = 45359.142
$ws.Range("M34").Value = -110463.996
$ws.Range("N34").Value = -45527.142
$ws.Range("H39").Value = 6538.0625
$ws.Range("J39").Value = 6873.8667
$ws.Range("L39").Value = 20621.6001
$ws.Range("N39").Value = -21209.6001
$ws.Range("H55").Value = 5247.5
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 5330
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 15990
$ws.Range("M55").Value = -14823
$ws.Range("N55").Value = -16344
$ws.Range("H109").Value = 5896.6665
$ws.Range("I109").Value = 1000
$ws.Range("J109").Value = 6876
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 20628
$ws.Range("M109").Value = -1960
$ws.Range("N109").Value = -22708
$ws.Range("H113").Value = 634.88135
$ws.Range("I113").Value = 638.85364
$ws.Range("J113").Value = 625.8333
$ws.Range("K113").Value = 1916.56092
$ws.Range("L113").Value = 1877.4999
$ws.Range("M113").Value = 253.4390799999999
$ws.Range("N113").Value = -6217.4999
$ws.Range("H135").Value = 1046.1628
$ws.Range("J135").Value = 2966.7273
$ws.Range("L135").Value = 26700.5457
$ws.Range("N135").Value = -31770.5457
$ws.Range("H138").Value = 2445
$ws.Range("I138").Value = 1238
$ws.Range("J138").Value = 4456.6665
$ws.Range("K138").Value = 3714
$ws.Range("L138").Value = 13369.9995
$ws.Range("M138").Value = 1426
$ws.Range("N138").Value = -23649.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 38500
$ws.Range("J52").Value = 38500
$ws.Range("L52").Value = 38500
$ws.Range("N52").Value = -39018
$ws.Range("H132").Value = 3442.3872
$ws.Range("I132").Value = 2585.9443
$ws.Range("J132").Value = 4628.231
$ws.Range("K132").Value = 7757.8329
$ws.Range("L132").Value = 13884.693
$ws.Range("M132").Value = -5227.8329
$ws.Range("N132").Value = -18944.693

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3489.4358
$ws.Range("I122").Value = 2902.7273
$ws.Range("J122").Value = 6716.3335
$ws.Range("K122").Value = 8708.1819
$ws.Range("L122").Value = 20149.0005
$ws.Range("M122").Value = -6258.1819
$ws.Range("N122").Value = -25049.0005
$ws.Range("H132").Value = 3979.2585
$ws.Range("I132").Value = 1478.7037
$ws.Range("J132").Value = 6157.161
$ws.Range("K132").Value = 4436.1111
$ws.Range("L132").Value = 18471.483
$ws.Range("M132").Value = -1906.1111
$ws.Range("N132").Value = -23531.483
$ws.Range("H136").Value = 2698.3635
$ws.Range("I136").Value = 1541
$ws.Range("J136").Value = 3755.087
$ws.Range("K136").Value = 4623
$ws.Range("L136").Value = 11265.261
$ws.Range("M136").Value = -2073
$ws.Range("N136").Value = -16365.261
$ws.Range("H141").Value = 41742.2
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 41742.2
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 41742.2
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -52102.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 28347.8
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 28347.8
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 28347.8
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -29367.8
$ws.Range("H113").Value = 8067.615
$ws.Range("I113").Value = 10339.9
$ws.Range("K113").Value = 31019.7
$ws.Range("M113").Value = -28849.7
$ws.Range("H132").Value = 6292871.5
$ws.Range("I132").Value = 5616.75
$ws.Range("J132").Value = 10103329
$ws.Range("K132").Value = 16850.25
$ws.Range("L132").Value = 30309987
$ws.Range("M132").Value = -14320.25
$ws.Range("N132").Value = -30315047
